$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds numeric-looking text (e.g. "1.000", "240.13").
# Force it to text format first so Excel stores these as text rather
# than coercing them into actual numbers, matching the source data.
$rngD = $ws.Range("D2:D51")
$rngD.NumberFormat = "@"

# Build a proper rectangular 2-D array (50 rows x 4 cols: B,C,D,E)
# for a single bulk Range.Value assignment.
$data = New-Object 'object[,]' 50,4
$data[0,0] = 'Bitcoin'; $data[0,1] = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; $data[0,2] = '29.337.45'; $data[0,3] = '  -0.10%  '
$data[1,0] = 'Ethereum'; $data[1,1] = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; $data[1,2] = '1.843.57'; $data[1,3] = '  -0.27%  '
$data[2,0] = 'TetherUSD'; $data[2,1] = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; $data[2,2] = '0.9994'; $data[2,3] = '  +0.05%  '
$data[3,0] = 'BNB'; $data[3,1] = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; $data[3,2] = '240.13'; $data[3,3] = '  -0.18%  '
$data[4,0] = 'XRP'; $data[4,1] = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; $data[4,2] = '0.6280'; $data[4,3] = '  -0.72%  '
$data[5,0] = 'USDC'; $data[5,1] = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; $data[5,2] = '1.000'; $data[5,3] = '  +0.10%  '
$data[6,0] = 'Dogecoin'; $data[6,1] = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; $data[6,2] = '0.07413'; $data[6,3] = '  -2.23%  '
$data[7,0] = 'Cardano'; $data[7,1] = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; $data[7,2] = '0.2893'; $data[7,3] = '  -1.08%  '
$data[8,0] = 'Solana'; $data[8,1] = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; $data[8,2] = '24.78'; $data[8,3] = '  +1.37%  '
$data[9,0] = 'TRON'; $data[9,1] = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; $data[9,2] = '0.07752'; $data[9,3] = '  +0.10%  '
$data[10,0] = 'WrappedEther'; $data[10,1] = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; $data[10,2] = '1.828.49'; $data[10,3] = '  -1.06%  '
$data[11,0] = 'Polkadot'; $data[11,1] = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; $data[11,2] = '4.990'; $data[11,3] = '  -0.62%  '
$data[12,0] = 'Polygon'; $data[12,1] = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; $data[12,2] = '0.6785'; $data[12,3] = '  -0.15%  '
$data[13,0] = 'ShibaInu'; $data[13,1] = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; $data[13,2] = '0.00001018'; $data[13,3] = '  -2.72%  '
$data[14,0] = 'Litecoin'; $data[14,1] = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; $data[14,2] = '82.09'; $data[14,3] = '  -1.33%  '
$data[15,0] = 'Uniswap'; $data[15,1] = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; $data[15,2] = '6.267'; $data[15,3] = '  +2.22%  '
$data[16,0] = 'WrappedBTC'; $data[16,1] = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; $data[16,2] = '29.357.97'; $data[16,3] = '  +0.00%  '
$data[17,0] = 'BitcoinCash'; $data[17,1] = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; $data[17,2] = '229.31'; $data[17,3] = '  -0.16%  '
$data[18,0] = 'Avalanche'; $data[18,1] = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; $data[18,2] = '12.31'; $data[18,3] = '  -0.37%  '
$data[19,0] = 'Dai'; $data[19,1] = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; $data[19,2] = '1.000'; $data[19,3] = '  +0.11%  '
$data[20,0] = 'Chainlink'; $data[20,1] = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; $data[20,2] = '7.442'; $data[20,3] = '  -0.23%  '
$data[21,0] = 'BinanceUSD'; $data[21,1] = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; $data[21,2] = '1.000'; $data[21,3] = '  +0.08%  '
$data[22,0] = 'Monero'; $data[22,1] = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; $data[22,2] = '158.69'; $data[22,3] = '  +0.06%  '
$data[23,0] = 'Cosmos'; $data[23,1] = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; $data[23,2] = '8.475'; $data[23,3] = '  +0.26%  '
$data[24,0] = 'Stellar'; $data[24,1] = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; $data[24,2] = '0.1354'; $data[24,3] = '  -2.92%  '
$data[25,0] = 'EthereumClassic'; $data[25,1] = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; $data[25,2] = '17.46'; $data[25,3] = '  -1.06%  '
$data[26,0] = 'Hedera'; $data[26,1] = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; $data[26,2] = '0.06636'; $data[26,3] = '  +16.76%  '
$data[27,0] = 'Toncoin'; $data[27,1] = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; $data[27,2] = '1.460'; $data[27,3] = '  +2.89%  '
$data[28,0] = 'PancakeSwap'; $data[28,1] = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; $data[28,2] = '1.486'; $data[28,3] = '  +0.88%  '
$data[29,0] = 'Filecoin'; $data[29,1] = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; $data[29,2] = '4.073'; $data[29,3] = '  -1.20%  '
$data[30,0] = 'InternetComputer(DFINITY)'; $data[30,1] = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; $data[30,2] = '4.071'; $data[30,3] = '  +0.62%  '
$data[31,0] = 'LidoDAOToken'; $data[31,1] = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; $data[31,2] = '1.837'; $data[31,3] = '  +0.51%  '
$data[32,0] = 'ARBITRUM'; $data[32,1] = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; $data[32,2] = '1.139'; $data[32,3] = '  -1.35%  '
$data[33,0] = 'ImmutableX'; $data[33,1] = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; $data[33,2] = '0.6934'; $data[33,3] = '  -1.97%  '
$data[34,0] = 'HuobiToken'; $data[34,1] = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; $data[34,2] = '2.575'; $data[34,3] = '  -0.16%  '
$data[35,0] = 'VeChain'; $data[35,1] = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; $data[35,2] = '0.01858'; $data[35,3] = '  +1.71%  '
$data[36,0] = 'MXToken'; $data[36,1] = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; $data[36,2] = '2.822'; $data[36,3] = '  +3.73%  '
$data[37,0] = 'Maker'; $data[37,1] = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; $data[37,2] = '1.246.08'; $data[37,3] = '  +0.14%  '
$data[38,0] = 'FraxShare'; $data[38,1] = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; $data[38,2] = '6.778'; $data[38,3] = '  +5.58%  '
$data[39,0] = 'TrustWalletToken'; $data[39,1] = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; $data[39,2] = '0.9372'; $data[39,3] = '  +3.84%  '
$data[40,0] = 'PaxDollar'; $data[40,1] = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; $data[40,2] = '1.000'; $data[40,3] = '  +0.11%  '
$data[41,0] = 'RocketPoolETH'; $data[41,1] = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'; $data[41,2] = '2.019.59'; $data[41,3] = '  +0.64%  '
$data[42,0] = 'Quant'; $data[42,1] = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; $data[42,2] = '100.77'; $data[42,3] = '  -0.82%  '
$data[43,0] = 'Aave'; $data[43,1] = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; $data[43,2] = '65.61'; $data[43,3] = '  -0.23%  '
$data[44,0] = 'Aptos'; $data[44,1] = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; $data[44,2] = '7.046'; $data[44,3] = '  -1.25%  '
$data[45,0] = 'RenderToken'; $data[45,1] = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; $data[45,2] = '1.709'; $data[45,3] = '  +2.06%  '
$data[46,0] = 'EnergySwap'; $data[46,1] = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; $data[46,2] = '9.019'; $data[46,3] = '  -0.16%  '
$data[47,0] = 'BabyDogeCoin'; $data[47,1] = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; $data[47,2] = '0.00000000116'; $data[47,3] = '  -0.05%  '
$data[48,0] = 'Algorand'; $data[48,1] = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; $data[48,2] = '0.1151'; $data[48,3] = '  -1.35%  '
$data[49,0] = 'TheSandbox'; $data[49,1] = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; $data[49,2] = '0.3902'; $data[49,3] = '  -1.33%  '

$ws.Range("B2:E51").Value = $data

# Clear the temporary text-number-format override so the price cells
# keep their original (default/unstyled) appearance.
$rngD.Style = "Normal"

Write-Output "applied cryptos update"
